$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
Write-Host $ws.Shapes.Item(5).TopLeftCell.Address()
